$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.876980721950531
$ws.Range("B1").Value = 2.671411275863647
$ws.Range("C1").Value = 4.680989265441895
$ws.Range("D1").Value = 2.209546327590942
$ws.Range("E1").Value = 1.303865313529968
